$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.310403
$ws.Range("H2").Value = 0.9312090000000001
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.162337
$ws.Range("N2").Value = 0.487011
$ws.Range("O2").Value = 0.1525168947008936
$ws.Range("P2").Value = 0.1525168947008936
$ws.Range("Q2").Value = 0.050389891811
$ws.Range("R2").Value = 0.453509026299
$ws.Range("S2").Value = 0.1525168947008936
$ws.Range("T2").Value = 0.1525168947008936

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.310403
$ws.Range("H3").Value = 0.9312090000000001
$ws.Range("O3").Value = 0.6673509415904804
$ws.Range("P3").Value = 0.6673509415904804
$ws.Range("Q3").Value = 0.2204853554923333
$ws.Range("R3").Value = 1.984368199431
$ws.Range("S3").Value = 0.6673509415904804
$ws.Range("T3").Value = 0.6673509415904804

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.310403
$ws.Range("H4").Value = 0.9312090000000001
$ws.Range("O4").Value = 0.180132163708626
$ws.Range("P4").Value = 0.180132163708626
$ws.Range("Q4").Value = 0.05951367065766668
$ws.Range("R4").Value = 0.5356230359190001
$ws.Range("S4").Value = 0.180132163708626
$ws.Range("T4").Value = 0.180132163708626
